$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.003.47'

$ws.Range('D3').Value = '3.182.67'

$ws.Range('E4').Value = '  -0.09%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '596.03'
$ws.Range('E5').Value = '  +3.71%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '153.54'
$ws.Range('E6').Value = '  +2.34%  '

$ws.Range('E7').Value = '  +0.02%  '

$ws.Range('D8').Value = '3.179.71'
$ws.Range('E8').Value = '  +0.59%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.538'
$ws.Range('E9').Value = '  +1.89%  '

$ws.Range('E10').Value = '  -0.88%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.07'
$ws.Range('E11').Value = '  -0.46%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.516'
$ws.Range('E12').Value = '  +3.55%  '

$ws.Range('E13').Value = '  +0.52%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '38.98'
$ws.Range('E14').Value = '  +4.84%  '

$ws.Range('D15').Value = '3.700.58'
$ws.Range('E15').Value = '  +0.58%  '

$ws.Range('B16').Value = 'Polkadot'
$ws.Range('C16').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '7.45'
$ws.Range('E16').Value = '  +4.97%  '

$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').Value = '66.019.65'
$ws.Range('E17').Value = '  +1.37%  '

$ws.Range('D18').Value = '3.176.49'
$ws.Range('E18').Value = '  +0.72%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.111'
$ws.Range('E19').Value = '  +0.57%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '509.00'
$ws.Range('E20').Value = '  +0.47%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '15.37'
$ws.Range('E21').Value = '  +3.47%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.738'
$ws.Range('E22').Value = '  +2.58%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.04'
$ws.Range('E23').Value = '  +3.90%  '

$ws.Range('E24').Value = '  -1.90%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '84.86'
$ws.Range('E25').Value = '  +0.54%  '

$ws.Range('E26').Value = '  -0.12%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.26'
$ws.Range('E27').Value = '  +3.78%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.99'
$ws.Range('E28').Value = '  +2.87%  '

$ws.Range('E29').Value = '  +5.40%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.06'
$ws.Range('E30').Value = '  +14.04%  '

$ws.Range('E31').Value = '  +2.47%  '

$ws.Range('E32').Value = '  +1.74%  '

$ws.Range('E33').Value = '  +2.26%  '

$ws.Range('E34').Value = '  +0.09%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.49'
$ws.Range('E35').Value = '  -0.62%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '54.77'
$ws.Range('E36').Value = '  -0.15%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0903'
$ws.Range('E37').Value = '  +0.15%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '483.23'
$ws.Range('E38').Value = '  +3.31%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0420'
$ws.Range('E39').Value = '  -0.09%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '8.83'
$ws.Range('E40').Value = '  +1.58%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.302'
$ws.Range('E41').Value = '  +6.74%  '

$ws.Range('E42').Value = '  +3.72%  '

$ws.Range('E43').Value = '  -5.36%  '

$ws.Range('D44').Value = '0.0₃0657'
$ws.Range('E44').Value = '  +11.64%  '

$ws.Range('D45').Value = '2.895.19'
$ws.Range('E45').Value = '  -4.99%  '

$ws.Range('E46').Value = '  -0.95%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '28.50'
$ws.Range('E47').Value = '  -0.15%  '

$ws.Range('E48').Value = '  +0.03%  '

$ws.Range('E49').Value = '  +1.60%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.31'
$ws.Range('E50').Value = '  +2.59%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.60'
$ws.Range('E51').Value = '  +1.51%  '
